# Auto-generated edit script: update cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Value)
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Value
    $r.NumberFormat = "General"
}

Set-TextValue "D2" "58.349.17"
Set-TextValue "E2" "  -3.72%  "

Set-TextValue "D3" "2.615.60"
Set-TextValue "E3" "  -3.11%  "

Set-TextValue "E4" "  -0.04%  "

Set-TextValue "D5" "520.29"
Set-TextValue "E5" "  -1.24%  "

Set-TextValue "D6" "142.94"
Set-TextValue "E6" "  -2.15%  "

Set-TextValue "E7" "  +0.28%  "

Set-TextValue "D8" "0.568"
Set-TextValue "E8" "  -1.62%  "

Set-TextValue "D9" "6.63"
Set-TextValue "E9" "  -2.56%  "

Set-TextValue "E10" "  -2.28%  "

Set-TextValue "E12" "  +1.11%  "

Set-TextValue "D13" "3.077.24"
Set-TextValue "E13" "  -3.13%  "

Set-TextValue "D14" "58.324.61"
Set-TextValue "E14" "  -3.80%  "

Set-TextValue "D15" "20.98"
Set-TextValue "E15" "  -1.78%  "

Set-TextValue "E16" "  -1.36%  "

Set-TextValue "D17" "2.606.64"
Set-TextValue "E17" "  -3.46%  "

$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D18" "336.47"
Set-TextValue "E18" "  -2.13%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D19" "4.40"
Set-TextValue "E19" "  -2.33%  "

Set-TextValue "E20" "  -1.44%  "

Set-TextValue "E21" "  -2.05%  "

Set-TextValue "E22" "  -0.08%  "

Set-TextValue "D23" "64.66"
Set-TextValue "E23" "  +1.69%  "

Set-TextValue "D24" "0.416"
Set-TextValue "E24" "  -0.92%  "

Set-TextValue "E25" "  -1.69%  "

Set-TextValue "E26" "  +0.26%  "

Set-TextValue "D27" "7.13"
Set-TextValue "E27" "  -1.89%  "

Set-TextValue "D28" "0.0₃0790"
Set-TextValue "E28" "  -3.58%  "

Set-TextValue "D29" "6.57"
Set-TextValue "E29" "  -3.32%  "

Set-TextValue "D30" "0.999"
Set-TextValue "E30" "  +0.08%  "

Set-TextValue "E31" "  -0.51%  "

Set-TextValue "E32" "  -1.13%  "

Set-TextValue "D33" "150.30"
Set-TextValue "E33" "  +0.43%  "

Set-TextValue "D34" "4.12"
Set-TextValue "E34" "  -3.45%  "

Set-TextValue "E35" "  -3.57%  "

Set-TextValue "E36" "  -4.87%  "

Set-TextValue "D37" "0.856"
Set-TextValue "E37" "  -2.70%  "

Set-TextValue "D38" "36.33"
Set-TextValue "E38" "  -2.21%  "

Set-TextValue "E39" "  -6.01%  "

Set-TextValue "E40" "  -0.69%  "

Set-TextValue "D41" "0.999"
Set-TextValue "E41" "  +0.39%  "

Set-TextValue "D42" "0.602"
Set-TextValue "E42" "  -1.37%  "

Set-TextValue "E43" "  -1.81%  "

Set-TextValue "D44" "269.48"
Set-TextValue "E44" "  -4.33%  "

Set-TextValue "D45" "10.64"
Set-TextValue "E45" "  +1.11%  "

Set-TextValue "D46" "19.17"
Set-TextValue "E46" "  -4.98%  "

Set-TextValue "D47" "0.0532"
Set-TextValue "E47" "  -1.46%  "

Set-TextValue "D48" "2.035.10"
Set-TextValue "E48" "  -4.53%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D49" "0.0229"
Set-TextValue "E49" "  -1.49%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D50" "4.68"
Set-TextValue "E50" "  -5.46%  "

Set-TextValue "D51" "18.30"
Set-TextValue "E51" "  -4.89%  "

Write-Output "Applied cryptos list update"
